# Update "想去人数" (column F) values on the "展览" and "全部类型" worksheets.
# Both worksheets contain identical data, and the same set of rows changed.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

# Mapping of row number -> new value for column F
$updates = @{
    3  = 3127
    5  = 166
    7  = 1732
    12 = 1420
    14 = 552
    15 = 355
    20 = 3
    23 = 115
    24 = 3341
    25 = 404
    26 = 208
    27 = 452
    28 = 29
    31 = 1056
    32 = 111
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
